$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure the new text cells use the same text number format as A2/B2 (style index 1)
$ws.Range("C2").NumberFormat = "@"
$ws.Range("F2").NumberFormat = "@"

# Fill in row 2 with the new GameServer data
$ws.Range("F2").Value = "127.0.0.1"
$ws.Range("A2").Value = "GameServer_1"
$ws.Range("B2").Value = "000104001"
$ws.Range("C2").Value = "GameServer_1"
$ws.Range("D2").Value = 5000
$ws.Range("E2").Value = 1
$ws.Range("G2").Value = 4001

# Update the data validation range so it no longer covers the now-filled F2 cell
$ws.Range("F2:F1048576").Validation.Delete()
$ws.Range("F3:F1048576").Validation.Add(3, 1, 1, '"TRUE,FALSE"')

# Update the active selection
$ws.Range("G3").Select()
